$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Column indices: 1 = Distance Class (m), 2 = N, 3 = Mantel r, 4 = p
# Row 1 is the header row; data rows are 2..16.
#
# Each entry: row, column, new text, optional new bold flag ($null = leave as-is)
$changes = @(
    @(2,  3, "0.093",  $null),
    @(3,  3, "-0.003", $null),
    @(3,  4, "0.404",  $null),
    @(4,  3, "-0.009", $null),
    @(4,  4, "0.505",  $false),
    @(5,  3, "-0.006", $null),
    @(5,  4, "0.758",  $false),
    @(6,  3, "0.001",  $null),
    @(6,  4, "1",      $null),
    @(7,  3, "-0.001", $null),
    @(7,  4, "1",      $null),
    @(8,  3, "-0.003", $null),
    @(8,  4, "1",      $null),
    @(9,  3, "-0.008", $null),
    @(9,  4, "1",      $null),
    @(10, 3, "-0.010", $null),
    @(10, 4, "1",      $null),
    @(11, 3, "-0.002", $null),
    @(12, 3, "-0.008", $null),
    @(12, 4, "1",      $null),
    @(13, 3, "-0.011", $null),
    @(13, 4, "1",      $null),
    @(14, 3, "0.003",  $null),
    @(14, 4, "1",      $null),
    @(15, 3, "-0.019", $null),
    @(15, 4, "0.377",  $null),
    @(16, 3, "-0.011", $null),
    @(16, 4, "1",      $false)
)

foreach ($chg in $changes) {
    $row = $chg[0]
    $col = $chg[1]
    $newText = $chg[2]
    $newBold = $chg[3]

    $cell = $t.Cell($row, $col)
    $r = $cell.Range
    # Trim the trailing cell-mark (\r\a) so we only touch the visible text.
    $r.MoveEnd(1, -1)
    $r.Select()
    $sel = $word.Selection

    if ($newBold -ne $null) {
        if ($newBold) {
            $sel.Font.Bold = 1
        } else {
            $sel.Font.Bold = 0
        }
    }

    $sel.Text = $newText
}
